# Update column G ("K") values on the active sheet to reflect the
# regenerated strikeout counts (K) instead of the prior "Strike#" values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 4
    4  = 2
    5  = 1
    6  = 0
    7  = 2
    8  = 1
    9  = 2
    10 = 2
    11 = 0
    12 = 3
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 4
    19 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
